# ECP-1197: adds date signed to import / export; also now defaults to
# date of upload instead of leaving the field empty.
#
# This inserts a new "Date Signed" column right after "Lease Amendment
# State" (i.e. before the old "Lease Reference" column), pushing all
# subsequent columns one place to the right, populates the new column
# for the rows that had a value, and moves the active selection along
# with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column C - this shifts the previous C:P into D:Q
# together with their values/styles/column widths intact.
$ws.Columns.Item(3).Insert()

# Give the new column roughly the same width as its neighbour (column B)
# - the sheet author never ran AutoFit/BestFit on this new column.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Header for the new column.
$ws.Cells.Item(1, 3).Value = "Date Signed"

# Populate "Date Signed" for the rows that have one, matching the date
# formatting (yyyy-mm-dd) already used elsewhere on the sheet - copy the
# format from a neighbouring date cell (now column E, "Start Date") so we
# reuse the existing style rather than create a new one.
$ws.Cells.Item(2, 5).Copy() | Out-Null
$ws.Cells.Item(2, 3).PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$ws.Cells.Item(2, 3).Value = 44004                        # 2020-06-22

$ws.Cells.Item(4, 5).Copy() | Out-Null
$ws.Cells.Item(4, 3).PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$ws.Cells.Item(4, 3).Value = 44002                        # 2020-06-20

$excel.CutCopyMode = 0

# Move the active selection from B4 to C4 (the new column), matching
# where the author was last working in the sheet.
$ws.Range("C4").Select() | Out-Null
